$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

$ws.Range("A$row").Value = "2025-08-17 04:06:05 UTC"
$ws.Range("B$row").Value = "2025-08-17 09:36:05 IST"
$ws.Range("C$row").Value = "SKIPPED"
$ws.Range("D$row").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E$row").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("F$row").Value = ""
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = ""

$ws.Range("A$row`:H$row").HorizontalAlignment = -4108
$ws.Range("A$row`:H$row").VerticalAlignment = -4108
